$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values (column G) for rows 2-23, replacing former Strike# values
$kValues = @{
    2  = 0
    3  = 0
    4  = 2
    5  = 0
    6  = 1
    7  = 4
    8  = 2
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 5
    14 = 2
    15 = 3
    16 = 2
    17 = 1
    18 = 0
    19 = 1
    20 = 1
    21 = 3
    22 = 0
    23 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
